# Weekly update: insert a new price record as row 30 (Fruta / hortaliza, semanal).
# This pushes the existing rows 30..124 down to 31..125, and the new row 30
# carries the same Mercado/Producto/Variedad/Calidad/Volumen/Unidad/Origen
# values as the (now-shifted) row below it, but with a newer Fecha and
# updated price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 30; everything below shifts down one row,
# inheriting the formatting (incl. the date number format on column D)
# from the row that used to be row 30.
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the new weekly record.
$ws.Cells.Item(30, 1).Value = 7
$ws.Cells.Item(30, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(30, 3).Value = "Ñuble"
$ws.Cells.Item(30, 4).Value = 44624
$ws.Cells.Item(30, 5).Value = 16
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100109
$ws.Cells.Item(30, 8).Value = "Uva"
$ws.Cells.Item(30, 9).Value = 100109001
$ws.Cells.Item(30, 10).Value = "Uva"
$ws.Cells.Item(30, 11).Value = "Red Globe"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 120
$ws.Cells.Item(30, 14).Value = 9000
$ws.Cells.Item(30, 15).Value = 10000
$ws.Cells.Item(30, 16).Value = 9500
$ws.Cells.Item(30, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(30, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 19).Value = 528
$ws.Cells.Item(30, 20).Value = 18
